$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (environment ssurgwsoadev4-oci...): update policy number and claim date
$ws.Range("E2").Value = 12112002243
$ws.Range("G2").Value = "'29/04/2022"

# Row 3 (environment i-preproducciongestion...): update policy number (as text) and claim date
$ws.Range("E3").Value = "'12112001841"
$ws.Range("G3").Value = "'07/04/2021"

# Update active selection to E3
$ws.Range("E3").Select() | Out-Null
